# Adding Registration Page Test Cases
# Adds a new "Sheet3" worksheet (after the existing Sheet2) containing a
# table of registration-form test data, mirroring what was pasted in from
# the author's testing notes.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.86
$ws.Columns.Item(3).ColumnWidth = 18.14

# --- Header / label rows ------------------------------------------------
$ws.Range("A1").Value = "Username:"
$ws.Range("A2").Value = "Password:"
$ws.Range("A3").Value = "Password confirmation:"

$ws.Range("A4").Value = "Username"
$ws.Range("B4").Value = "Password"
$ws.Range("C4").Value = "Confirm Password"

# --- Test case data rows -------------------------------------------------
$ws.Range("A5").Value = "ninja tester"
$ws.Range("B5").Value = "ninja@123"
$ws.Range("C5").Value = "ninja@123"

$ws.Range("A6").Value = "ninjatesterss"
$ws.Range("B6").Value = "nin@1"
$ws.Range("C6").Value = "nin@1"

$ws.Range("A7").Value = "ninjatesterss"
$ws.Range("B7").Value = 1234567890
$ws.Range("C7").Value = 1234567890

$ws.Range("A8").Value = "ninjatesterss"
$ws.Range("B8").Value = "ninja@123"
$ws.Range("C8").Value = "ninja@12345"

$ws.Range("A9").Value = "ninjatesterss"
$ws.Range("B9").Value = "ninja@123"
$ws.Range("C9").Value = "ninja@123"

$ws.Range("A10").Value = "ninjatestersssss"
$ws.Range("B10").Value = "ninja@123"
$ws.Range("C10").Value = "ninja@123"

$ws.Range("A11").Value = "numpyninjatester1"
$ws.Range("B11").Value = "ninja@123"
$ws.Range("C11").Value = "ninja@123"

$ws.Range("A12").Value = "numpyninjatester2"
$ws.Range("B12").Value = "ninja@123"
$ws.Range("C12").Value = "ninja@123"

$ws.Range("A13").Value = "numpyninjatester3"
$ws.Range("B13").Value = "ninja@123"
$ws.Range("C13").Value = "ninja@123"

# --- Formatting ------------------------------------------------------------
# Whole used range gets the plain, readingOrder-aware default font used
# elsewhere in the workbook.
$used = $ws.Range("A1:C13")
$used.Font.Name = "Calibri"

# A5 ("ninja tester") carries its own small Arial style, left aligned.
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 8
$a5.Font.Color = 2039583
$a5.HorizontalAlignment = -4131
